# Add a "seq_run_id" column (with a barcode-mismatch / run-id parameter
# value) into the sample sheet template, inserted right after the
# sample_id column (new column D), shifting the existing D:M columns to
# E:N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D; this shifts the old D:M -> E:N
# and carries the existing formatting of column D along (Excel's normal
# "Insert" behaviour), which we override right after for the new column.
$ws.Range("D1").EntireColumn.Insert()

# Fill in the header + the constant run-id value for every data row.
$ws.Range("D1").Value = "seq_run_id"
$ws.Range("D2:D13").Value = "210810_NB551189_0075_AHVMWGBGXJ"

# The new column is sized explicitly (not auto "best fit" like its
# neighbours). Target stored width is 38.1640625 characters; the COM
# width setter here quantizes to whole pixels (width = (round(6*w)+5)/6),
# so 37.33 is the input that lands closest (38.1666...) to that value.
$ws.Range("D1").ColumnWidth = 37.33

# Match the saved selection/active cell from the authored edit.
$ws.Range("D16").Select()
